$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.122.84"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'2.055.75"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'250.07"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'0.669"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'60.24"
$ws.Range("E7").Value = "  +9.19%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.388"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").Value = "'0.0794"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "'16.10"
$ws.Range("E12").Value = "  +7.84%  "
$ws.Range("D13").Value = "'2.355.27"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'0.830"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").Value = "'5.73"
$ws.Range("E15").Value = "  +9.48%  "
$ws.Range("D16").Value = "'2.063.78"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "'18.55"
$ws.Range("E17").Value = "  +30.39%  "
$ws.Range("D18").Value = "'37.161.33"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'75.76"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").Value = "'0.0₃0908"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").Value = "'5.47"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").Value = "'239.13"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  +13.42%  "
$ws.Range("E26").Value = "  +6.42%  "
$ws.Range("D27").Value = "'169.16"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'20.29"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'1.14"
$ws.Range("E30").Value = "  +9.61%  "
$ws.Range("D31").Value = "'4.84"
$ws.Range("E31").Value = "  +5.53%  "
$ws.Range("D32").Value = "'0.0627"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "'4.55"
$ws.Range("E33").Value = "  +4.54%  "
$ws.Range("D34").Value = "'0.0898"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  +5.17%  "
$ws.Range("D39").Value = "'1.35"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'5.17"
$ws.Range("E40").Value = "  +26.78%  "
$ws.Range("D41").Value = "'3.09"
$ws.Range("E41").Value = "  +10.35%  "
$ws.Range("D42").Value = "'17.86"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "'97.75"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").Value = "'3.92"
$ws.Range("E47").Value = "  -5.44%  "
$ws.Range("D48").Value = "'1.293.43"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'6.88"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'2.238.79"
$ws.Range("E51").Value = "  -0.58%  "